$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36 ---
$ws.Cells.Item(36, 1).Value = "'14343571"
$ws.Cells.Item(36, 1).Style = "Normal"

$ws.Cells.Item(36, 2).Value = "'2025-08-04"
$ws.Cells.Item(36, 2).Style = "Normal"

$ws.Cells.Item(36, 3).Value = "Radu Albot"

$ws.Cells.Item(36, 4).Value = "Tadeas Paroulek"

$ws.Cells.Item(36, 5).Value = "Gana Tadeas Paroulek"

$ws.Cells.Item(36, 6).Value = 2.5

$ws.Cells.Item(36, 7).Value = "'"
$ws.Cells.Item(36, 7).Style = "Normal"

$ws.Cells.Item(36, 8).Value = "'"
$ws.Cells.Item(36, 8).Style = "Normal"

# --- Row 37 ---
$ws.Cells.Item(37, 1).Value = "'14344407"
$ws.Cells.Item(37, 1).Style = "Normal"

$ws.Cells.Item(37, 2).Value = "'2025-08-04"
$ws.Cells.Item(37, 2).Style = "Normal"

$ws.Cells.Item(37, 3).Value = "Zdenek Kolar"

$ws.Cells.Item(37, 4).Value = "Eric Vanshelboim"

$ws.Cells.Item(37, 5).Value = "Gana Eric Vanshelboim"

$ws.Cells.Item(37, 6).Value = 3

$ws.Cells.Item(37, 7).Value = "'"
$ws.Cells.Item(37, 7).Style = "Normal"

$ws.Cells.Item(37, 8).Value = "'"
$ws.Cells.Item(37, 8).Style = "Normal"
